$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute(" #Report-Date#")
Write-Output "found: $found start=$($r.Start) end=$($r.End)"
$r.Text = " #Repor Date Here#"
Write-Output "set done"
